$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 2 (shifts existing rows 2-10 down to 3-11,
# carrying their formatting/styles with them).
$ws.Rows(2).Insert()

# Populate the newly inserted row 2 with the new award entry.
$ws.Range("A2").Value = "IX Excellence Awards"
$ws.Range("B2").Value = "Dic. 2022"
$ws.Range("C2").Value = "Universidad El Bosque"
$ws.Range("D2").Value = "Bogota, Colombia"
$ws.Range("E2").Value = "COP`$10.000.000"

# Match the style already used by the equivalent "amount" column elsewhere
# in the table (row 4's E cell, same as rows 5/6).
$ws.Range("E2").Style = $ws.Range("E4").Style

# Update the view to match the edited state (topLeftCell / selection).
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("B3:B4").Select()
